# save data done + era data updated
# Adds a new "Save" column (H) mirroring the existing header style,
# with a value of 0 for the single data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last header cell (G1) onto the new header
# cell (H1) so the new column matches the existing bold/centered/bordered
# header style, then set its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# New data cell for the "Save" column.
$ws.Range("H2").Value = 0
